# Aggiunti attributi al file Auto
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 4 new rows before row 3, pushing the existing content
#    (old row3..row7) down to row7..row11
$ws.Range("A3:A6").EntireRow.Insert() | Out-Null

# 2. Fill in the new legend rows (3 and 4). Column J3 stays empty.
$ws.Range("B3").Value = "NOTNULL"
$ws.Range("C3:I3").Value = "NOTNULL"
$ws.Range("K3:L3").Value = "NOTNULL"
$ws.Range("N3").Value = "NOTNULL"

$ws.Range("A3").Value = "IDENTITY, UNIQUE"

$ws.Range("D4").Value = "UNIQUE"

$ws.Range("M3").Value = "DEFAULT"

# 3. Apply the new (italic) styling to the legend block.
#    Build each combined style once on a scratch cell, then copy the
#    resulting format onto the target range - this avoids generating
#    extra unused intermediate cell styles.
$scratchB = $ws.Range("Z100")
$scratchB.HorizontalAlignment = -4108
$scratchB.Font.Italic = $true
$scratchB.Copy() | Out-Null
$ws.Range("B3:N5").PasteSpecial(-4122) | Out-Null
$scratchB.Clear() | Out-Null

$scratchA = $ws.Range("Z101")
$scratchA.HorizontalAlignment = -4108
$scratchA.VerticalAlignment = -4108
$scratchA.WrapText = $true
$scratchA.Font.Italic = $true
$scratchA.Copy() | Out-Null
$ws.Range("A3:A5").PasteSpecial(-4122) | Out-Null
$scratchA.Clear() | Out-Null

$excel.CutCopyMode = 0

# 4. Merge the legend label cell across the 3 new rows.
$ws.Range("A3:A5").Merge() | Out-Null

# 5. Resize some columns.
$ws.Columns("A").ColumnWidth = 9.5
$ws.Columns("K").ColumnWidth = 10.833333333333334
$ws.Columns("L").ColumnWidth = 10.5

# 6. Update the active selection.
$ws.Range("C4").Select() | Out-Null

Write-Host "Done"
